$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows right before the current row 310, shifting
# all subsequent rows down by two (old row 310 -> new row 312, etc.)
$ws.Rows.Item(310).Insert()
$ws.Rows.Item(310).Insert()

# New row 310
$ws.Range("A310").Value = 7
$ws.Range("B310").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C310").Value = "Ñuble"
$ws.Range("D310").Value = 44985
$ws.Range("E310").Value = 16
$ws.Range("F310").Value = 100112023
$ws.Range("G310").Value = "Brócoli"
$ws.Range("H310").Value = "Sin especificar"
$ws.Range("I310").Value = "Primera"
$ws.Range("J310").Value = 400
$ws.Range("K310").Value = 850
$ws.Range("L310").Value = 900
$ws.Range("M310").Value = 875
$ws.Range("N310").Value = "$/unidad"
$ws.Range("O310").Value = "Región del Maule"
$ws.Range("P310").Value = 875
$ws.Range("Q310").Value = 1
$ws.Range("R310").Value = "Hortaliza"

# New row 311
$ws.Range("A311").Value = 7
$ws.Range("B311").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C311").Value = "Ñuble"
$ws.Range("D311").Value = 44985
$ws.Range("E311").Value = 16
$ws.Range("F311").Value = 100112023
$ws.Range("G311").Value = "Brócoli"
$ws.Range("H311").Value = "Sin especificar"
$ws.Range("I311").Value = "Segunda"
$ws.Range("J311").Value = 300
$ws.Range("K311").Value = 700
$ws.Range("L311").Value = 750
$ws.Range("M311").Value = 725
$ws.Range("N311").Value = "$/unidad"
$ws.Range("O311").Value = "Región del Maule"
$ws.Range("P311").Value = 725
$ws.Range("Q311").Value = 1
$ws.Range("R311").Value = "Hortaliza"
